# Automatic update of files.
#
# The workbook's "Förändrad" (last-changed) column (C) is bumped by one day
# for every data row, from serial 45202 (2023-10-03) to serial 45203
# (2023-10-04). This touches every data row, C2:C180.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C180").Value = 45203
